$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Logistic Regression -> Logistic Regression (Tuned)
$ws.Range("A2").Value = "Logistic Regression (Tuned)"
$ws.Range("B2").Value = 0.9755065433601691
$ws.Range("C2").Value = 0.5714285714285714
$ws.Range("D2").Value = 0.8235294117647058

# Row 3: Random Forest -> Random Forest (Tuned)
$ws.Range("A3").Value = "Random Forest (Tuned)"
$ws.Range("B3").Value = 0.9652347204357263
$ws.Range("C3").Value = 0.7959183673469388
$ws.Range("D3").Value = 0.8764044943820225

# Row 4: Gradient Boosting -> Gradient Boosting (Tuned)
$ws.Range("A4").Value = "Gradient Boosting (Tuned)"
$ws.Range("B4").Value = 0.9129154739414056
$ws.Range("C4").Value = 0.7448979591836735
$ws.Range("D4").Value = 0.7604166666666666

# Row 5: XGBoost -> XGBoost (Tuned)
$ws.Range("A5").Value = "XGBoost (Tuned)"
$ws.Range("B5").Value = 0.9758012852721281
$ws.Range("C5").Value = 0.6938775510204082
$ws.Range("D5").Value = 0.8947368421052632

# Row 6: LightGBM -> LightGBM (Tuned)
$ws.Range("A6").Value = "LightGBM (Tuned)"
$ws.Range("B6").Value = 0.5479759978696036
$ws.Range("C6").Value = 0.173469387755102
$ws.Range("D6").Value = 0.08292682926829269

$wb.Save()
